# Doing Updates for Financials
# Insert a new fiscal-year column (FY2018, period ending 2018-12-31) as the
# new first data column on the CBFV sheet. All the existing yearly data in
# columns D:K shifts right to E:L, and the new column D is populated with
# the latest year's figures for the Income Statement, Balance Sheet and
# Cash Flow Statement sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before D; this shifts D:K -> E:L (formats move
#    with the cells, including the date format on row 7/38/80 and the
#    #,##0 number format used throughout the data rows).
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D cells default to an unformatted style,
#    so copy the number formatting from column E (which just received the
#    original column D's formatting) back onto column D. Only touch the
#    three data blocks (Income Statement, Balance Sheet, Cash Flow
#    Statement) so we don't introduce stray styled cells on blank/label
#    rows that never had a D:K value to begin with.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate new column D with the FY2018 figures.

# Income Statement
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 43600
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -1500
$ws.Range("D17").Value = 8500
$ws.Range("D18").Value = 35200
$ws.Range("D20").Value = -26600
$ws.Range("D21").Value = 11600
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 8600
$ws.Range("D24").Value = 1500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 7100
$ws.Range("D27").Value = 7100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 26600
$ws.Range("D33").Value = 7100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 7100

# Balance Sheet
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 16600
$ws.Range("D42").Value = 36700
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 23400
$ws.Range("D49").Value = 39400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 4900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1281300
$ws.Range("D57").Value = 6000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1143700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 57800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 137600
$ws.Range("D77").Value = 0

# Cash Flow Statement
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 7100
$ws.Range("D83").Value = 3000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 13700
$ws.Range("D91").Value = -4400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 27500
$ws.Range("D96").Value = -4500
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -8400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 32700
